$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5000
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("M64").Value = -4752
$ws.Range("H67").Value = 5000
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("M67").Value = -4142
$ws.Range("H116").Value = 5286.2666
$ws.Range("I116").Value = 4468.636
$ws.Range("K116").Value = 4468.636
$ws.Range("M116").Value = -1026.636
$ws.Range("H132").Value = 5900.15
$ws.Range("I132").Value = 5765.8237
$ws.Range("K132").Value = 17297.4711
$ws.Range("M132").Value = -14767.4711
$ws.Range("H137").Value = 2842.2
$ws.Range("I137").Value = 2223.3333
$ws.Range("K137").Value = 6669.999899999999
$ws.Range("M137").Value = -4119.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").ClearContents()
$ws.Range("N61").Value = 0
$ws.Range("H63").Value = 3502.5
$ws.Range("I63").Value = 3005
$ws.Range("K63").Value = 3005
$ws.Range("M63").Value = -2319
$ws.Range("H66").Value = 3502.5
$ws.Range("I66").Value = 3005
$ws.Range("K66").Value = 15025
$ws.Range("M66").Value = -11593
$ws.Range("H74").Value = 1653.6923
$ws.Range("J74").Value = 2325
$ws.Range("L74").Value = 2325
$ws.Range("N74").Value = -4073
$ws.Range("H77").Value = 1653.6923
$ws.Range("J77").Value = 2325
$ws.Range("L77").Value = 11625
$ws.Range("N77").Value = -20361
$ws.Range("H88").Value = 7284.4287
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 7284.4287
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = 7284.4287
$ws.Range("N88").Value = -8096.4287
$ws.Range("H91").Value = 7284.4287
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 7284.4287
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = 7284.4287
$ws.Range("N91").Value = -10092.4287
$ws.Range("H132").Value = 3859.4736
$ws.Range("I132").Value = 3430.353
$ws.Range("J132").Value = 7507
$ws.Range("K132").Value = 10291.059
$ws.Range("L132").Value = 22521
$ws.Range("M132").Value = -7761.059000000001
$ws.Range("N132").Value = -27581
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").ClearContents()
$ws.Range("N136").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5950
$ws.Range("I99").Value = 5950
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 5950
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -4452
$ws.Range("H105").Value = 3013.1333
$ws.Range("I105").Value = 2749.5715
$ws.Range("K105").Value = 2749.5715
$ws.Range("M105").Value = -1002.5715
$ws.Range("H127").Value = 25998
$ws.Range("J127").Value = 25998
$ws.Range("L127").Value = 25998
$ws.Range("N127").Value = -35918

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 9573.375
$ws.Range("J62").Value = 9624.25
$ws.Range("L62").Value = 9624.25
$ws.Range("N62").Value = -10872.25
$ws.Range("H65").Value = 9573.375
$ws.Range("J65").Value = 9624.25
$ws.Range("L65").Value = 48121.25
$ws.Range("N65").Value = -54361.25
$ws.Range("H99").Value = 8100
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H109").Value = 56128.5
$ws.Range("I109").Value = 53259
$ws.Range("J109").Value = 58998
$ws.Range("K109").Value = 53259
$ws.Range("L109").Value = 58998
$ws.Range("M109").Value = -52219
$ws.Range("N109").Value = -61078
$ws.Range("H126").Value = 8100
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 27826.285
$ws.Range("I14").Value = 27826.285
$ws.Range("K14").Value = 83478.855
$ws.Range("M14").Value = -83305.855
$ws.Range("H137").Value = 6842.875
$ws.Range("J137").Value = 7549
$ws.Range("L137").Value = 22647
$ws.Range("N137").Value = -32847

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = 0
$ws.Range("H70").Value = 7571
$ws.Range("I70").Value = 6749.5
$ws.Range("K70").Value = 6749.5
$ws.Range("M70").Value = -6479.5
$ws.Range("H73").Value = 7571
$ws.Range("I73").Value = 6749.5
$ws.Range("K73").Value = 6749.5
$ws.Range("M73").Value = -5813.5
$ws.Range("H80").Value = 4993
$ws.Range("I80").Value = 2852.5
$ws.Range("J80").Value = 5849.2
$ws.Range("K80").Value = 2852.5
$ws.Range("L80").Value = 5849.2
$ws.Range("M80").Value = -1854.5
$ws.Range("N80").Value = -7845.2
$ws.Range("H83").Value = 4993
$ws.Range("I83").Value = 2852.5
$ws.Range("J83").Value = 5849.2
$ws.Range("K83").Value = 14262.5
$ws.Range("L83").Value = 29246
$ws.Range("M83").Value = -9270.5
$ws.Range("N83").Value = -39230
$ws.Range("H132").Value = 736.8
$ws.Range("I132").Value = 785.44446
$ws.Range("J132").Value = 299
$ws.Range("K132").Value = 2356.33338
$ws.Range("L132").Value = 897
$ws.Range("M132").Value = 173.66662
$ws.Range("N132").Value = -5957

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 119999
$ws.Range("J109").Value = 119999
$ws.Range("L109").Value = 119999
$ws.Range("N109").Value = -122773
$ws.Range("H113").Value = 1167.6316
$ws.Range("I113").Value = 777.8182
$ws.Range("J113").Value = 1703.625
$ws.Range("K113").Value = 2333.4546
$ws.Range("L113").Value = 5110.875
$ws.Range("M113").Value = -163.4546
$ws.Range("N113").Value = -9450.875
